# Excel COM-interop script implementing the commit:
# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# Functional changes applied to Hoja1:
#   1. Update "VALOR MORA" total (E11): 341640 -> 512460
#   2. Update "Cant. Periodos" (F13): 2 -> 3
#   3. Add a third period (2509) block of 3 worker rows to the table,
#      matching the existing pattern used for periods 2507 / 2508,
#      pushing the signature footer block down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update header totals -------------------------------------------------
$ws.Range("E11").Value2 = 512460
$ws.Range("F13").Value2 = 3

# --- 2. Insert 3 new rows for the new period (2509) --------------------------
# Row 21 currently holds the last worker of period 2508 and carries the
# special "closing" border style used only by the final row of the table.
# Insert 3 blank rows right after it; the table footer (signature block)
# is automatically pushed down by Excel's row-insert shifting.
$ws.Range("A22:A24").EntireRow.Insert()

# Preserve the "closing" border formatting (currently still on row 21) by
# copying it onto the new last row (24) before row 21 is reformatted.
$ws.Range("B21:J21").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)   # xlPasteFormats

# Row 21 is no longer the last row of the table, so give it the regular
# (non-closing) formatting used by the rest of the data rows, matching
# row 20's style. Apply the same regular formatting to the two other new
# rows (22-23).
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B22:J23").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3. Fill in the new period-2509 rows, mirroring the existing pattern -----
$ws.Range("B22").Value2 = "CC"
$ws.Range("C22").Value2 = "1044912108"
$ws.Range("D22").Value2 = "JOSE ULISES CASTRO OSORIO"
$ws.Range("E22").Value2 = "2509"
$ws.Range("F22").Value2 = 56940
$ws.Range("G22").Value2 = 1423500

$ws.Range("B23").Value2 = "CC"
$ws.Range("C23").Value2 = "73352461"
$ws.Range("D23").Value2 = "HERNANDO ENRIQUE DE ARCO BOSSIO"
$ws.Range("E23").Value2 = "2509"
$ws.Range("F23").Value2 = 56940
$ws.Range("G23").Value2 = 1423500

$ws.Range("B24").Value2 = "CC"
$ws.Range("C24").Value2 = "85167757"
$ws.Range("D24").Value2 = "EUSEBIO MANUEL FERREIRA ALFARO"
$ws.Range("E24").Value2 = "2509"
$ws.Range("F24").Value2 = 56940
$ws.Range("G24").Value2 = 1423500

$wb.Save()
